$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A195").Value = "2023-12-11 16:19:31"
$ws.Range("B195").Value = 0.0014

$ws.Range("A196").Value = "2023-12-11 16:19:38"
$ws.Range("B196").Value = 0.0004

$ws.Range("A197").Value = "2023-12-11 16:19:51"
$ws.Range("B197").Value = 0.0004
